$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: rotate the daily price records across rows 2-7
# (Fecha, Calidad, Volumen, Precio minimo, Precio maximo,
#  Precio promedio ponderado, Precio $/Kg).

$updates = @{
    2 = @{ D = 44162; L = "Primera"; M = 100; N = 7000; O = 7000; P = 7000; S = 4667 }
    3 = @{ D = 44162; L = "Segunda"; M = 100; N = 6500; O = 6500; P = 6500; S = 4333 }
    4 = @{ D = 44176; L = "Primera"; M = 300; N = 5000; O = 6000; P = 5500; S = 3667 }
    5 = @{ D = 44169; L = "Primera"; M = 400; N = 5500; O = 6000; P = 5750; S = 3833 }
    6 = @{ D = 44159; L = "Segunda"; M = 200; N = 6500; O = 7000; P = 6750; S = 4500 }
    7 = @{ D = 44166; L = "Primera"; M = 200; N = 6000; O = 6500; P = 6250; S = 4167 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
